# Team03Report.xlsx - "Working on Sprint 2" commit
#
# Applies the real content edits described by the diff:
#  - Backlog: assign Sprint numbers to newly-started/finished stories and
#    flip a few stories from "in progress" to "done"
#  - Burndown: append the next sprint data point (and matching chart range)
#  - Sprint2: mark the first few stories done, fill in their actual
#    size/time/completion-date columns, and add the "Review Results"
#    retro section (mirroring the one already on Sprint1)
#  - leave the workbook looking at the Stories sheet, matching where the
#    author ended up navigating to while doing this work

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Backlog: sprint assignment + status updates
# ---------------------------------------------------------------------
$backlog = $wb.Worksheets.Item("Backlog")

# Stories that moved into Sprint 2 and were completed
$backlog.Cells.Item(4, 1).Value = 2
$backlog.Cells.Item(4, 5).Value = "done"

$backlog.Cells.Item(5, 1).Value = 2
$backlog.Cells.Item(5, 5).Value = "done"

$backlog.Cells.Item(8, 1).Value = 2
$backlog.Cells.Item(8, 5).Value = "done"

# Story that moved into Sprint 2 but is still in progress
$backlog.Cells.Item(9, 1).Value = 2

# Stories picked up for Sprint 3
$backlog.Cells.Item(10, 1).Value = 3
$backlog.Cells.Item(10, 4).Value = "zd"
$backlog.Cells.Item(10, 5).Value = "in progress"

$backlog.Cells.Item(16, 1).Value = 3
$backlog.Cells.Item(16, 4).Value = "zd"
$backlog.Cells.Item(16, 5).Value = "in progress"

# ---------------------------------------------------------------------
# Burndown: new data point for the sprint boundary + burndown chart
# ---------------------------------------------------------------------
$burndown = $wb.Worksheets.Item("Burndown")
$burndown.Cells.Item(4, 1).Value = 40982
$burndown.Cells.Item(4, 2).Value = 8
$burndown.Cells.Item(4, 3).Value = 4
$burndown.Cells.Item(4, 4).Value = 800
$burndown.Cells.Item(4, 5).Value = 390
$burndown.Cells.Item(4, 6).Value = 30

# ---------------------------------------------------------------------
# Sprint2: completed stories get actual size/time + completion date;
# add the Review Results retro section like Sprint1 already has
# ---------------------------------------------------------------------
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint2 = $wb.Worksheets.Item("Sprint2")

$sprint2.Cells.Item(2, 4).Value = "done"
$sprint2.Cells.Item(2, 7).Value = "100 LOC"
$sprint2.Cells.Item(2, 8).Value = "1 hour"
$sprint2.Cells.Item(2, 9).Value = 40982

$sprint2.Cells.Item(3, 4).Value = "done"
$sprint2.Cells.Item(3, 7).Value = "100 LOC"
$sprint2.Cells.Item(3, 8).Value = "1 hour"
$sprint2.Cells.Item(3, 9).Value = 40982

$sprint2.Cells.Item(4, 4).Value = "done"
$sprint2.Cells.Item(4, 7).Value = "50 LOC"
$sprint2.Cells.Item(4, 8).Value = "30 mins"
$sprint2.Cells.Item(4, 9).Value = 40982

$sprint2.Cells.Item(5, 9).Value = 40982

# Match the date formatting already used in the "Act Time" column on Sprint1
$sprint1.Range("I2").Copy()
$sprint2.Range("I2:I5").PasteSpecial(-4122)

# Copy the whole "Review Results" retro block from Sprint1 onto Sprint2
$reviewSrc = $sprint1.Range("B14:B22")
$reviewDst = $sprint2.Range("B9:B17")
$reviewSrc.Copy()
$reviewDst.PasteSpecial(-4163)
$reviewSrc.Copy()
$reviewDst.PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Navigation: move through the sheets and finish up on Stories
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Sprint1").Activate()
$wb.Worksheets.Item("Sprint1").Range("I10").Select()

$backlog.Activate()
$backlog.Range("E16").Select()

$burndown.Activate()
$burndown.Range("J18").Select()

$sprint2.Activate()
$sprint2.Range("E14").Select()

$wb.Worksheets.Item("Sprint3").Activate()
$wb.Worksheets.Item("Sprint3").Range("F32").Select()

$stories = $wb.Worksheets.Item("Stories")
$stories.Activate()
$stories.Range("C9").Select()
